# Add the two new "date fields" rows (value / budget) to the "details" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("details")

$ws.Range("B17").Value = "value"
$ws.Range("C17").Value = 4000

$ws.Range("B18").Value = "budget"
$ws.Range("C18").Value = 300

$ws.Range("C17").Select() | Out-Null
